$d = $word.ActiveDocument

# Replace each unique text value (date header + multiplication equations)
# using Find/Execute with exact matching against the whole document content.
$d.Content.Find.Execute("2025-09-02 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-03 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("442×9=3978", $true, $false, $false, $false, $false, $true, 1, $false, "520×5=2600", 2) | Out-Null
$d.Content.Find.Execute("271×2=542", $true, $false, $false, $false, $false, $true, 1, $false, "249×6=1494", 2) | Out-Null
$d.Content.Find.Execute("769×5=3845", $true, $false, $false, $false, $false, $true, 1, $false, "442×7=3094", 2) | Out-Null
$d.Content.Find.Execute("368×7=2576", $true, $false, $false, $false, $false, $true, 1, $false, "377×5=1885", 2) | Out-Null
$d.Content.Find.Execute("743×6=4458", $true, $false, $false, $false, $false, $true, 1, $false, "382×6=2292", 2) | Out-Null
$d.Content.Find.Execute("271×5=1355", $true, $false, $false, $false, $false, $true, 1, $false, "468×9=4212", 2) | Out-Null
$d.Content.Find.Execute("573×4=2292", $true, $false, $false, $false, $false, $true, 1, $false, "610×3=1830", 2) | Out-Null
$d.Content.Find.Execute("240×8=1920", $true, $false, $false, $false, $false, $true, 1, $false, "374×9=3366", 2) | Out-Null
$d.Content.Find.Execute("326×8=2608", $true, $false, $false, $false, $false, $true, 1, $false, "171×9=1539", 2) | Out-Null
$d.Content.Find.Execute("343×9=3087", $true, $false, $false, $false, $false, $true, 1, $false, "966×9=8694", 2) | Out-Null
$d.Content.Find.Execute("584×8=4672", $true, $false, $false, $false, $false, $true, 1, $false, "942×9=8478", 2) | Out-Null
$d.Content.Find.Execute("235×9=2115", $true, $false, $false, $false, $false, $true, 1, $false, "710×6=4260", 2) | Out-Null
$d.Content.Find.Execute("487×8=3896", $true, $false, $false, $false, $false, $true, 1, $false, "135×6=810", 2) | Out-Null
$d.Content.Find.Execute("872×2=1744", $true, $false, $false, $false, $false, $true, 1, $false, "576×4=2304", 2) | Out-Null
$d.Content.Find.Execute("181×6=1086", $true, $false, $false, $false, $false, $true, 1, $false, "429×8=3432", 2) | Out-Null
$d.Content.Find.Execute("296×9=2664", $true, $false, $false, $false, $false, $true, 1, $false, "383×3=1149", 2) | Out-Null
$d.Content.Find.Execute("226×4=904", $true, $false, $false, $false, $false, $true, 1, $false, "502×9=4518", 2) | Out-Null
$d.Content.Find.Execute("913×8=7304", $true, $false, $false, $false, $false, $true, 1, $false, "918×3=2754", 2) | Out-Null
$d.Content.Find.Execute("489×2=978", $true, $false, $false, $false, $false, $true, 1, $false, "821×5=4105", 2) | Out-Null
$d.Content.Find.Execute("305×4=1220", $true, $false, $false, $false, $false, $true, 1, $false, "906×8=7248", 2) | Out-Null
$d.Content.Find.Execute("736×2=1472", $true, $false, $false, $false, $false, $true, 1, $false, "798×4=3192", 2) | Out-Null
$d.Content.Find.Execute("271×3=813", $true, $false, $false, $false, $false, $true, 1, $false, "843×6=5058", 2) | Out-Null
$d.Content.Find.Execute("117×6=702", $true, $false, $false, $false, $false, $true, 1, $false, "473×2=946", 2) | Out-Null
$d.Content.Find.Execute("934×3=2802", $true, $false, $false, $false, $false, $true, 1, $false, "606×9=5454", 2) | Out-Null
$d.Content.Find.Execute("704×5=3520", $true, $false, $false, $false, $false, $true, 1, $false, "294×8=2352", 2) | Out-Null

Write-Host "Replacements complete."
